# Apply the changes described by the diff:
#  1. Rename the three "Include from SNOMED CT[ N]" sheets to "Include #0/1/2".
#  2. On the Metadata sheet:
#       - update the "Date" value
#       - insert a new "Jurisdiction" row (with an empty value) right after "Contact",
#         pushing Description/Purpose/Copyright/Immutable down by one row.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from SNOMED CT" sheets -------------------------
$wb.Worksheets.Item("Include from SNOMED CT").Name   = "Include #0"
$wb.Worksheets.Item("Include from SNOMED CT 2").Name = "Include #1"
$wb.Worksheets.Item("Include from SNOMED CT 3").Name = "Include #2"

# --- 2. Update the Metadata sheet ------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date property value (row 8, column B)
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new row above row 11 ("Description") for the "Jurisdiction" property
$ws.Rows.Item(11).Insert()

# Copy formatting from the row above (Contact, row 10) onto the newly inserted row
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" row with its label and an empty value
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
